# Auto-generated edit script applying scheduled market-data refresh
# to Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 830.1111
$ws.Range("I18").Value = 853
$ws.Range("J18").Value = 750
$ws.Range("K18").Value = 853
$ws.Range("L18").Value = 750
$ws.Range("M18").Value = -569
$ws.Range("N18").Value = -1318

$ws.Range("H43").Value = 1532.6
$ws.Range("I43").Value = 1921.4286
$ws.Range("J43").Value = 1323.2307
$ws.Range("K43").Value = 1921.4286
$ws.Range("L43").Value = 1323.2307
$ws.Range("M43").Value = -1852.4286
$ws.Range("N43").Value = -1461.2307

$ws.Range("H98").Value = 569.35297
$ws.Range("I98").Value = 569.35297
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 569.35297
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 928.64703

$ws.Range("H122").Value = 569.35297
$ws.Range("I122").Value = 569.35297
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1708.05891
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 741.9410899999998

$ws.Range("H129").Value = 1300.6364
$ws.Range("I129").Value = 898.5
$ws.Range("J129").Value = 1390
$ws.Range("K129").Value = 2695.5
$ws.Range("L129").Value = 4170
$ws.Range("M129").Value = 2304.5
$ws.Range("N129").Value = -14170

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 50002
$ws.Range("I6").Value = 50002
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 50002
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -49829
$ws.Range("N6").ClearContents()

$ws.Range("H9").Value = 40000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 40000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 40000
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -40340

$ws.Range("H20").Value = 40000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 40000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 40000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -40540

$ws.Range("H63").Value = 5000
$ws.Range("I63").Value = 5000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 5000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -4314

$ws.Range("H66").Value = 5000
$ws.Range("I66").Value = 5000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 25000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -21568

$ws.Range("H132").Value = 10595.52
$ws.Range("I132").Value = 8152.2354
$ws.Range("J132").Value = 15787.5
$ws.Range("K132").Value = 24456.7062
$ws.Range("L132").Value = 47362.5
$ws.Range("M132").Value = -21926.7062
$ws.Range("N132").Value = -52422.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13350.857
$ws.Range("I82").Value = 7242.6665
$ws.Range("J82").Value = 50000
$ws.Range("K82").Value = 7242.6665
$ws.Range("L82").Value = 50000
$ws.Range("M82").Value = -6859.6665
$ws.Range("N82").Value = -50766

$ws.Range("H85").Value = 13350.857
$ws.Range("I85").Value = 7242.6665
$ws.Range("J85").Value = 50000
$ws.Range("K85").Value = 7242.6665
$ws.Range("L85").Value = 50000
$ws.Range("M85").Value = -5916.6665
$ws.Range("N85").Value = -52652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1880.0714
$ws.Range("I58").Value = 890.7
$ws.Range("J58").Value = 4353.5
$ws.Range("K58").Value = 890.7
$ws.Range("L58").Value = 4353.5
$ws.Range("M58").Value = -687.7
$ws.Range("N58").Value = -4759.5

$ws.Range("H122").Value = 1452.4
$ws.Range("I122").Value = 1302
$ws.Range("J122").Value = 1602.8
$ws.Range("K122").Value = 3906
$ws.Range("L122").Value = 4808.4
$ws.Range("M122").Value = -1456
$ws.Range("N122").Value = -9708.4

$ws.Range("H132").Value = 38467760
$ws.Range("I132").Value = 58831224
$ws.Range("J132").Value = 3439.7778
$ws.Range("K132").Value = 176493672
$ws.Range("L132").Value = 10319.3334
$ws.Range("M132").Value = -176491142
$ws.Range("N132").Value = -15379.3334

$ws.Range("H134").Value = 1242.7858
$ws.Range("I134").Value = 948.75
$ws.Range("J134").Value = 3007
$ws.Range("K134").Value = 2846.25
$ws.Range("L134").Value = 9021
$ws.Range("M134").Value = -311.25
$ws.Range("N134").Value = -14091

$ws.Range("H136").Value = 1880.0714
$ws.Range("I136").Value = 890.7
$ws.Range("J136").Value = 4353.5
$ws.Range("K136").Value = 2672.1
$ws.Range("L136").Value = 13060.5
$ws.Range("M136").Value = -122.1000000000004
$ws.Range("N136").Value = -18160.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1078.1666
$ws.Range("I122").Value = 492.5
$ws.Range("J122").Value = 2249.5
$ws.Range("K122").Value = 4432.5
$ws.Range("L122").Value = 20245.5
$ws.Range("M122").Value = -1982.5
$ws.Range("N122").Value = -25145.5

$ws.Range("H132").Value = 874.75
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 949.5
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 8545.5
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -13605.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1074.8334
$ws.Range("I126").Value = 1012.25
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 3036.75
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -566.75
$ws.Range("N126").Value = -8540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1341.3684
$ws.Range("I7").Value = 1193.5385
$ws.Range("J7").Value = 1661.6666
$ws.Range("K7").Value = 1193.5385
$ws.Range("L7").Value = 1661.6666
$ws.Range("M7").Value = -1081.5385
$ws.Range("N7").Value = -1885.6666

$ws.Range("H108").Value = 48000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 48000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 48000
$ws.Range("N108").Value = -55680

$ws.Range("H126").Value = 1341.3684
$ws.Range("I126").Value = 1193.5385
$ws.Range("J126").Value = 1661.6666
$ws.Range("K126").Value = 3580.6155
$ws.Range("L126").Value = 4984.9998
$ws.Range("M126").Value = -1110.6155
$ws.Range("N126").Value = -9924.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 44995
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 44995
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 44995
$ws.Range("N46").Value = -45457

$ws.Range("H126").Value = 612.1177
$ws.Range("I126").Value = 476.25
$ws.Range("J126").Value = 732.8889
$ws.Range("K126").Value = 1428.75
$ws.Range("L126").Value = 2198.6667
$ws.Range("M126").Value = 1041.25
$ws.Range("N126").Value = -7138.6667

$ws.Range("H132").Value = 23816874
$ws.Range("I132").Value = 41676796
$ws.Range("J132").Value = 3646
$ws.Range("K132").Value = 125030388
$ws.Range("L132").Value = 10938
$ws.Range("M132").Value = -125027858
$ws.Range("N132").Value = -15998

$ws.Range("H134").Value = 44995
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 44995
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 134985
$ws.Range("N134").Value = -140055
